$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "mayor/kuidadó no ta kustumá", $true, $false, $false, $false, $false,
    $true, 1, $false, "mayor/dunadónan di kuido no ta kustumá", 2)

$d.Content.Find.Execute(
    "adolesentenan, mayornan/kuidadónan tambe", $true, $false, $false, $false, $false,
    $true, 1, $false, "adolesentenan, mayornan/dunadónan di kuido tambe", 2)

$d.Content.Find.Execute(
    "kontestá na nan yunan", $true, $false, $false, $false, $false,
    $true, 1, $false, "kontestá nan yunan", 2)

$d.Content.Find.Execute(
    "kustumbrá ku ta puntra", $true, $false, $false, $false, $false,
    $true, 1, $false, "kustumá ku ta puntra", 2)

$d.Content.Find.Execute(
    "skuchá nan, mester", $true, $false, $false, $false, $false,
    $true, 1, $false, "skucha nan, mester", 2)
